$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 and 33: Monero / NEARProtocol swapped positions with refreshed data
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "172.89"
$ws.Range("E32").Value = "  +1.80%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  +10.32%  "

# Refreshed price / volume figures for the remaining coins
$ws.Range("D2").Value = "64.203.95"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "2.784.42"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("D5").Value = "587.78"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.30"
$ws.Range("E6").Value = "  +7.87%  "
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "3.279.17"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").Value = "27.55"
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("D15").Value = "64.111.94"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  +6.02%  "
$ws.Range("D17").Value = "2.796.91"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "5.08"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").Value = "368.04"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "7.08"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").Value = "0.578"
$ws.Range("E22").Value = "  +8.24%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "67.63"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").Value = "8.89"
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("E27").Value = "  +13.35%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("E31").Value = "  +6.12%  "
$ws.Range("D34").Value = "20.91"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  +5.73%  "
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "4.29"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.30"
$ws.Range("E40").Value = "  +12.02%  "
$ws.Range("D41").Value = "341.71"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("D42").Value = "40.31"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").Value = "22.55"
$ws.Range("E43").Value = "  +4.37%  "
$ws.Range("D44").Value = "22.59"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("D45").Value = "0.0611"
$ws.Range("E45").Value = "  +3.30%  "
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").Value = "138.77"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").Value = "2.175.41"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("E51").Value = "  +0.32%  "
